$wb = $excel.ActiveWorkbook

# Correct the 1in diameters in column P (rows 1-21) from 0.127 to 0.0127
$ws1 = $wb.Worksheets.Item("1in")
for ($r = 1; $r -le 21; $r++) {
    $ws1.Cells.Item($r, 16).Value = 0.0127
}

# Select P1:P21 with active cell P1 on the "1in" sheet, and make it the active tab
$ws1.Activate()
$ws1.Range("P1:P21").Select()
